{"js": "// Resume edit: \"Get the PCB on there\"\n//\n// 1. \"Replatformed Java application onto Docker container\" paragraph:\n//    split the run so \"Replatformed\" is wrapped in spell-check proofErr\n//    markers (spellStart/spellEnd) just like Word does after a spell check\n//    pass, without changing the visible text.\n// 2. \"Designing an oscilloscope app for audio signals for use in \" paragraph:\n//    reworded to \"Developing an audio-based oscilloscope app\" (split across\n//    4 runs).\n// 3. \"Advanced Embedded Systems class\" paragraph: turned into its own\n//    bullet (numId 12) reading \"Designing a Line-out to Mic-in adapter PCB\".\n// 4. \"Internet of Things \u2013 C++ \\tFall 2021\" paragraph: the trailing\n//    \" \" run and \"Fall 2021\" run are merged into a single \" Fall 2021\" run.\n\nconst RFONTS = '<w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>';\n\nfunction flatOpc(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' + bodyXml + '<w:sectPr/></w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Replace a whole paragraph's contents (pPr + runs) via a crafted OOXML\n// fragment. This gives exact control over run boundaries / proofErr\n// placement that the plain text APIs (which silently coalesce runs with\n// identical formatting) cannot express.\nfunction replaceParagraph(paragraph, innerXml) {\n  const range = paragraph.getRange(\"Whole\");\n  range.insertOoxml(flatOpc(\"<w:p>\" + innerXml + \"</w:p>\"), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idxReplatformed = -1;\nlet idxOscilloscope = -1;\nlet idxEmbeddedSystems = -1;\nlet idxIot = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Replatformed Java application onto Docker container \") {\n    idxReplatformed = i;\n  } else if (t === \"Designing an oscilloscope app for audio signals for use in \") {\n    idxOscilloscope = i;\n  } else if (t === \"Advanced Embedded Systems class\") {\n    idxEmbeddedSystems = i;\n  } else if (t === \"Internet of Things \u2013 C++\\t Fall 2021\") {\n    idxIot = i;\n  }\n}\n\nif (idxReplatformed === -1 || idxOscilloscope === -1 || idxEmbeddedSystems === -1 || idxIot === -1) {\n  throw new Error(\n    \"Could not locate all target paragraphs: \" +\n      JSON.stringify({ idxReplatformed, idxOscilloscope, idxEmbeddedSystems, idxIot })\n  );\n}\n\n// 1. \"Replatformed Java application onto Docker container \"\n{\n  const pPr =\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + RFONTS + '</w:rPr></w:pPr>';\n  const inner =\n    pPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t>Replatformed</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t xml:space=\"preserve\"> Java application onto Docker container </w:t></w:r>';\n  replaceParagraph(paragraphs.items[idxReplatformed], inner);\n}\n\n// 2. \"Designing an oscilloscope app for audio signals for use in \" ->\n//    \"Developing an audio-based oscilloscope app\"\n{\n  const pPr =\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"12\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + RFONTS + '</w:rPr></w:pPr>';\n  const inner =\n    pPr +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t>Developing</w:t></w:r>' +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t xml:space=\"preserve\"> an </w:t></w:r>' +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t xml:space=\"preserve\">audio-based </w:t></w:r>' +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t>oscilloscope app</w:t></w:r>';\n  replaceParagraph(paragraphs.items[idxOscilloscope], inner);\n}\n\n// 3. \"Advanced Embedded Systems class\" -> \"Designing a Line-out to Mic-in\n//    adapter PCB\" and promote it to a bullet (numId 12), matching its\n//    sibling above.\n{\n  const pPr =\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"12\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + RFONTS + '</w:rPr></w:pPr>';\n  const inner =\n    pPr +\n    '<w:r><w:rPr>' + RFONTS + '</w:rPr><w:t>Designing a Line-out to Mic-in adapter PCB</w:t></w:r>';\n  replaceParagraph(paragraphs.items[idxEmbeddedSystems], inner);\n}\n\n// 4. \"Internet of Things \u2013 C++ \\t Fall 2021\": merge the trailing \" \" run\n//    and the \"Fall 2021\" run into a single \" Fall 2021\" run.\n{\n  const pPr =\n    '<w:pPr><w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + RFONTS + '</w:rPr></w:pPr>';\n  const inner =\n    pPr +\n    '<w:r><w:rPr>' + RFONTS + '<w:b/></w:rPr><w:t xml:space=\"preserve\">Internet of Things \u2013 </w:t></w:r>' +\n    '<w:r><w:rPr>' + RFONTS + '<w:bCs/></w:rPr><w:t>C++</w:t></w:r>' +\n    '<w:r w:rsidRPr=\"003B19FB\"><w:rPr>' + RFONTS + '</w:rPr><w:tab/><w:t xml:space=\"preserve\"> Fall 2021</w:t></w:r>';\n  replaceParagraph(paragraphs.items[idxIot], inner);\n}\n\nawait context.sync();\n", "ps1": "# Resume edit: \"Get the PCB on there\"\n#\n# 1. \"Replatformed Java application onto Docker container\" paragraph:\n#    split the run so \"Replatformed\" is wrapped in spell-check proofErr\n#    markers (spellStart/spellEnd) just like Word does after a spell check\n#    pass, without changing the visible text.\n# 2. \"Designing an oscilloscope app for audio signals for use in \" paragraph:\n#    reworded to \"Developing an audio-based oscilloscope app\" (split across\n#    4 runs).\n# 3. \"Advanced Embedded Systems class\" paragraph: turned into its own\n#    bullet (numId 12) reading \"Designing a Line-out to Mic-in adapter PCB\".\n# 4. \"Internet of Things - C++ \\tFall 2021\" paragraph: the trailing\n#    \" \" run and \"Fall 2021\" run are merged into a single \" Fall 2021\" run.\n\n$d = $word.ActiveDocument\n\n$RFONTS = '<w:rFonts w:asciiTheme=\"minorHAnsi\" w:hAnsiTheme=\"minorHAnsi\" w:cstheme=\"minorHAnsi\"/>'\n\nfunction Wrap-FlatOpc($bodyXml) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyXml + '<w:sectPr/></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# Replace a whole paragraph's contents (pPr + runs) via a crafted OOXML\n# fragment inserted through Range.InsertXML. This gives exact control over\n# run boundaries / proofErr placement that plain text assignment (which\n# silently coalesces runs with identical formatting) cannot express.\nfunction Replace-Paragraph($paragraph, $innerXml) {\n    $range = $paragraph.Range\n    [void]$range.InsertXML((Wrap-FlatOpc('<w:p>' + $innerXml + '</w:p>')))\n}\n\n# Locate the target paragraphs by their (trimmed) visible text so the\n# script is resilient to any earlier edits shifting paragraph indices.\n$idxReplatformed = -1\n$idxOscilloscope = -1\n$idxEmbeddedSystems = -1\n$idxIot = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Replatformed Java application onto Docker container \") {\n        $idxReplatformed = $i\n    } elseif ($t -eq \"Designing an oscilloscope app for audio signals for use in \") {\n        $idxOscilloscope = $i\n    } elseif ($t -eq \"Advanced Embedded Systems class\") {\n        $idxEmbeddedSystems = $i\n    } elseif ($t -eq \"Internet of Things \u2013 C++`t Fall 2021\") {\n        $idxIot = $i\n    }\n}\n\nif ($idxReplatformed -eq -1 -or $idxOscilloscope -eq -1 -or $idxEmbeddedSystems -eq -1 -or $idxIot -eq -1) {\n    throw \"Could not locate all target paragraphs ($idxReplatformed, $idxOscilloscope, $idxEmbeddedSystems, $idxIot)\"\n}\n\n# 1. \"Replatformed Java application onto Docker container \"\n$pPr1 = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"9\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'\n$inner1 = $pPr1 +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Replatformed</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space=\"preserve\"> Java application onto Docker container </w:t></w:r>'\nReplace-Paragraph $d.Paragraphs.Item($idxReplatformed) $inner1\n\n# 2. \"Designing an oscilloscope app for audio signals for use in \" ->\n#    \"Developing an audio-based oscilloscope app\"\n$pPr2 = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"12\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'\n$inner2 = $pPr2 +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Developing</w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space=\"preserve\"> an </w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t xml:space=\"preserve\">audio-based </w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>oscilloscope app</w:t></w:r>'\nReplace-Paragraph $d.Paragraphs.Item($idxOscilloscope) $inner2\n\n# 3. \"Advanced Embedded Systems class\" -> \"Designing a Line-out to Mic-in\n#    adapter PCB\" and promote it to a bullet (numId 12), matching its\n#    sibling above.\n$pPr3 = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"12\"/></w:numPr>' +\n    '<w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'\n$inner3 = $pPr3 +\n    '<w:r><w:rPr>' + $RFONTS + '</w:rPr><w:t>Designing a Line-out to Mic-in adapter PCB</w:t></w:r>'\nReplace-Paragraph $d.Paragraphs.Item($idxEmbeddedSystems) $inner3\n\n# 4. \"Internet of Things - C++ \\t Fall 2021\": merge the trailing \" \" run\n#    and the \"Fall 2021\" run into a single \" Fall 2021\" run.\n$pPr4 = '<w:pPr><w:tabs><w:tab w:val=\"left\" w:pos=\"720\"/><w:tab w:val=\"right\" w:pos=\"8640\"/></w:tabs>' +\n    '<w:rPr>' + $RFONTS + '</w:rPr></w:pPr>'\n$inner4 = $pPr4 +\n    '<w:r><w:rPr>' + $RFONTS + '<w:b/></w:rPr><w:t xml:space=\"preserve\">Internet of Things \u2013 </w:t></w:r>' +\n    '<w:r><w:rPr>' + $RFONTS + '<w:bCs/></w:rPr><w:t>C++</w:t></w:r>' +\n    '<w:r w:rsidRPr=\"003B19FB\"><w:rPr>' + $RFONTS + '</w:rPr><w:tab/><w:t xml:space=\"preserve\"> Fall 2021</w:t></w:r>'\nReplace-Paragraph $d.Paragraphs.Item($idxIot) $inner4\n"}
